$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Odds updates for Jogos da Semana FlashScore 2025-06-04 (row 3,4,5,6,8,9,13,14,16)
$ws.Range("G3").Value = 2.9
$ws.Range("H3").Value = 3.55
$ws.Range("I3").Value = 2.25
$ws.Range("T3").Value = 9.5
$ws.Range("U3").Value = 16.5
$ws.Range("V3").Value = 11.25
$ws.Range("W3").Value = 37
$ws.Range("X3").Value = 25
$ws.Range("Y3").Value = 35
$ws.Range("AE3").Value = 8.25
$ws.Range("AF3").Value = 12
$ws.Range("AG3").Value = 9.5
$ws.Range("AH3").Value = 24
$ws.Range("AI3").Value = 18.5
$ws.Range("J4").Value = 1.1
$ws.Range("K4").Value = 7
$ws.Range("G5").Value = 2
$ws.Range("I5").Value = 4.2
$ws.Range("G6").Value = 1.62
$ws.Range("H6").Value = 3.7
$ws.Range("I6").Value = 4.9
$ws.Range("L6").Value = 1.26
$ws.Range("M6").Value = 3.15
$ws.Range("N6").Value = 1.78
$ws.Range("O6").Value = 1.83
$ws.Range("P6").Value = 1.39
$ws.Range("Q6").Value = 2.57
$ws.Range("R6").Value = 1.78
$ws.Range("S6").Value = 1.82
$ws.Range("T6").Value = 6.8
$ws.Range("U6").Value = 7.6
$ws.Range("V6").Value = 8
$ws.Range("W6").Value = 12
$ws.Range("X6").Value = 13
$ws.Range("Y6").Value = 26
$ws.Range("Z6").Value = 10.75
$ws.Range("AA6").Value = 7.3
$ws.Range("AB6").Value = 16
$ws.Range("AC6").Value = 75
$ws.Range("AD6").Value = 600
$ws.Range("AE6").Value = 13.5
$ws.Range("AF6").Value = 29
$ws.Range("AG6").Value = 16
$ws.Range("AI6").Value = 50
$ws.Range("AJ6").Value = 50
$ws.Range("G8").Value = 2.67
$ws.Range("H8").Value = 3.1
$ws.Range("I8").Value = 2.55
$ws.Range("L8").Value = 1.35
$ws.Range("M8").Value = 2.72
$ws.Range("N8").Value = 2.02
$ws.Range("O8").Value = 1.62
$ws.Range("R8").Value = 1.78
$ws.Range("S8").Value = 1.83
$ws.Range("U8").Value = 13
$ws.Range("W8").Value = 30
$ws.Range("Z8").Value = 8.25
$ws.Range("AA8").Value = 6
$ws.Range("AD8").Value = 600
$ws.Range("AE8").Value = 7.6
$ws.Range("AG8").Value = 9.75
$ws.Range("AH8").Value = 28
$ws.Range("AI8").Value = 23
$ws.Range("AJ8").Value = 35
$ws.Range("S9").Value = 2.67
$ws.Range("AA9").Value = 9.25
$ws.Range("AF9").Value = 23
$ws.Range("G13").Value = 3.35
$ws.Range("I13").Value = 2.1
$ws.Range("M13").Value = 2.42
$ws.Range("N13").Value = 2.25
$ws.Range("R13").Value = 2.02
$ws.Range("T13").Value = 7.9
$ws.Range("U13").Value = 16
$ws.Range("V13").Value = 12.5
$ws.Range("W13").Value = 45
$ws.Range("X13").Value = 37
$ws.Range("Y13").Value = 55
$ws.Range("AB13").Value = 19
$ws.Range("AE13").Value = 5.9
$ws.Range("AF13").Value = 8.75
$ws.Range("AH13").Value = 18.5
$ws.Range("AI13").Value = 20
$ws.Range("G14").Value = 3.4
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 2.1
$ws.Range("Q14").Value = 2.2
$ws.Range("T14").Value = 7.6
$ws.Range("U14").Value = 16
$ws.Range("V14").Value = 13
$ws.Range("AA14").Value = 6.2
$ws.Range("AE14").Value = 5.6
$ws.Range("AF14").Value = 8.5
$ws.Range("AH14").Value = 19
$ws.Range("H16").Value = 3.45
$ws.Range("T16").Value = 10.75
$ws.Range("U16").Value = 15
$ws.Range("AE16").Value = 10.25
